$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" contain identical data tables and both need
# their F-column (想去人数 / "want to go" count) values bumped.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F4").Value = 1564
    $ws.Range("F7").Value = 398
    $ws.Range("F10").Value = 421
}
